$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2480.2
$ws.Range("I2").Value = 3000.25
$ws.Range("K2").Value = 3000.25
$ws.Range("M2").Value = -2887.25

$ws.Range("H16").Value = 8000
$ws.Range("I16").Value = 8000
$ws.Range("K16").Value = 8000
$ws.Range("M16").Value = -7770

$ws.Range("H33").Value = 303.9091
$ws.Range("I33").Value = 317.22223
$ws.Range("J33").Value = 244
$ws.Range("K33").Value = 317.22223
$ws.Range("L33").Value = 244
$ws.Range("M33").Value = -88.22223000000002
$ws.Range("N33").Value = -702

$ws.Range("H40").Value = 4800.75
$ws.Range("I40").Value = 3625.75
$ws.Range("J40").Value = 5975.75
$ws.Range("K40").Value = 3625.75
$ws.Range("L40").Value = 5975.75
$ws.Range("M40").Value = -3450.75
$ws.Range("N40").Value = -6325.75

$ws.Range("H45").Value = 7124.75
$ws.Range("J45").Value = 7124.75
$ws.Range("L45").Value = 21374.25
$ws.Range("N45").Value = -21758.25

$ws.Range("H86").Value = 2452
$ws.Range("J86").Value = 1904
$ws.Range("L86").Value = 1904
$ws.Range("N86").Value = -4150

$ws.Range("H89").Value = 2452
$ws.Range("J89").Value = 1904
$ws.Range("L89").Value = 9520
$ws.Range("N89").Value = -20752

$ws.Range("H98").Value = 7026.5557
$ws.Range("I98").Value = 8654.714
$ws.Range("K98").Value = 8654.714
$ws.Range("M98").Value = -7156.714

$ws.Range("H106").Value = 4111.625
$ws.Range("I106").Value = 4748.25
$ws.Range("K106").Value = 4748.25
$ws.Range("M106").Value = -4117.25

$ws.Range("H122").Value = 7026.5557
$ws.Range("I122").Value = 8654.714
$ws.Range("K122").Value = 25964.142
$ws.Range("M122").Value = -23514.142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7164.567
$ws.Range("I61").Value = 5404
$ws.Range("K61").Value = 5404
$ws.Range("M61").Value = -5192

$ws.Range("H105").Value = 97252.98
$ws.Range("J105").Value = 97252.98
$ws.Range("L105").Value = 97252.98
$ws.Range("N105").Value = -104240.98

$ws.Range("H122").Value = 4638
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4638
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").Value = 13914
$ws.Range("N122").Value = -18814

$ws.Range("H136").Value = 7164.567
$ws.Range("I136").Value = 5404
$ws.Range("K136").Value = 16212
$ws.Range("M136").Value = -13662

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1348.3
$ws.Range("I64").Value = 930.5
$ws.Range("J64").Value = 1975
$ws.Range("K64").Value = 930.5
$ws.Range("L64").Value = 1975
$ws.Range("M64").Value = -705.5
$ws.Range("N64").Value = -2425

$ws.Range("H67").Value = 1348.3
$ws.Range("I67").Value = 930.5
$ws.Range("J67").Value = 1975
$ws.Range("K67").Value = 930.5
$ws.Range("L67").Value = 1975
$ws.Range("M67").Value = -150.5
$ws.Range("N67").Value = -3535

$ws.Range("H80").Value = 95173.37
$ws.Range("J80").Value = 173734.67
$ws.Range("L80").Value = 173734.67
$ws.Range("N80").Value = -175730.67

$ws.Range("H83").Value = 95173.37
$ws.Range("J83").Value = 173734.67
$ws.Range("L83").Value = 868673.3500000001
$ws.Range("N83").Value = -878657.3500000001

$ws.Range("H86").Value = 62196.637
$ws.Range("I86").Value = 1584.2307
$ws.Range("K86").Value = 1584.2307
$ws.Range("M86").Value = -461.2307000000001

$ws.Range("H89").Value = 62196.637
$ws.Range("I89").Value = 1584.2307
$ws.Range("K89").Value = 7921.1535
$ws.Range("M89").Value = -2305.1535

$ws.Range("H105").Value = 27039470
$ws.Range("I105").Value = 47636050
$ws.Range("J105").Value = 6465.125
$ws.Range("K105").Value = 47636050
$ws.Range("L105").Value = 6465.125
$ws.Range("M105").Value = -47634303
$ws.Range("N105").Value = -9959.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2499
$ws.Range("J94").Value = 891.3570999999999
$ws.Range("L94").Value = 891.3570999999999
$ws.Range("N94").Value = -1793.3571

$ws.Range("H122").Value = 4044.889
$ws.Range("I122").Value = 3398.3333
$ws.Range("K122").Value = 10194.9999
$ws.Range("M122").Value = -7744.999899999999

$ws.Range("H132").Value = 27267.262
$ws.Range("I132").Value = 2816.611
$ws.Range("K132").Value = 8449.832999999999
$ws.Range("M132").Value = -5919.832999999999

$ws.Range("H134").Value = 4649.6875
$ws.Range("I134").Value = 2742
$ws.Range("K134").Value = 8226
$ws.Range("M134").Value = -5691

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1554197
$ws.Range("I4").Value = 895727
$ws.Range("K4").Value = 2687181
$ws.Range("M4").Value = -2687069

$ws.Range("H37").Value = 90965.55499999999
$ws.Range("J37").Value = 90965.55499999999
$ws.Range("L37").Value = 272896.665
$ws.Range("N37").Value = -273120.665

$ws.Range("H68").Value = 1315.1666
$ws.Range("J68").Value = 947.5
$ws.Range("L68").Value = 2842.5
$ws.Range("N68").Value = -4464.5

$ws.Range("H71").Value = 1315.1666
$ws.Range("J71").Value = 947.5
$ws.Range("L71").Value = 8527.5
$ws.Range("N71").Value = -16639.5

$ws.Range("H109").Value = 1974.25
$ws.Range("I109").Value = 1799.3334
$ws.Range("K109").Value = 5398.0002
$ws.Range("M109").Value = -4358.0002

$ws.Range("I131").Value = 47619860
$ws.Range("J131").Value = 9125.177
$ws.Range("K131").Value = 142859580
$ws.Range("L131").Value = 27375.531
$ws.Range("M131").Value = -142854540
$ws.Range("N131").Value = -37455.531

$ws.Range("H132").Value = 67994.664
$ws.Range("I132").Value = 77919.30499999999
$ws.Range("J132").Value = 3484.5
$ws.Range("K132").Value = 701273.7449999999
$ws.Range("L132").Value = 31360.5
$ws.Range("M132").Value = -698743.7449999999
$ws.Range("N132").Value = -36420.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 2858.3333
$ws.Range("J6").Value = 3348.4
$ws.Range("L6").Value = 3348.4
$ws.Range("N6").Value = -3574.4

$ws.Range("H14").Value = 192531.06
$ws.Range("I14").Value = 362344.22
$ws.Range("K14").Value = 362344.22
$ws.Range("M14").Value = -362176.22

$ws.Range("H16").Value = 2858.3333
$ws.Range("J16").Value = 3348.4
$ws.Range("L16").Value = 3348.4
$ws.Range("N16").Value = -3848.4

$ws.Range("H113").Value = 14251802
$ws.Range("I113").Value = 2054.5
$ws.Range("J113").Value = 28501550
$ws.Range("K113").Value = 2054.5
$ws.Range("L113").Value = 28501550
$ws.Range("M113").Value = 115.5
$ws.Range("N113").Value = -28505890

$ws.Range("H134").Value = 199999
$ws.Range("J134").Value = 199999
$ws.Range("L134").Value = 599997
$ws.Range("N134").Value = -605067

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 14954.546
$ws.Range("I20").Value = 14500
$ws.Range("K20").Value = 14500
$ws.Range("M20").Value = -14274

$ws.Range("H55").Value = 195.41176
$ws.Range("I55").Value = 149.8
$ws.Range("J55").Value = 214.41667
$ws.Range("K55").Value = 149.8
$ws.Range("L55").Value = 214.41667
$ws.Range("M55").Value = 23.19999999999999
$ws.Range("N55").Value = -560.4166700000001

$ws.Range("H106").Value = 18735.125
$ws.Range("J106").Value = 18735.125
$ws.Range("L106").Value = 18735.125
$ws.Range("N106").Value = -21259.125

$ws.Range("H122").Value = 5091.0527
$ws.Range("I122").Value = 4825
$ws.Range("J122").Value = 5547.143
$ws.Range("K122").Value = 14475
$ws.Range("L122").Value = 16641.429
$ws.Range("M122").Value = -12025
$ws.Range("N122").Value = -21541.429

$ws.Range("H130").Value = 14500
$ws.Range("J130").Value = 14500
$ws.Range("L130").Value = 14500
$ws.Range("N130").Value = -24540
